$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is an unambiguous number (e.g. 332.52) must be
# pre-formatted as Text so Excel stores them as strings (matching the original
# inline-string cells) instead of auto-converting them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.236.82"
$ws.Range("E2").Value = "  -2.89%  "
$ws.Range("D3").Value = "1.933.05"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("D5").Value = "332.52"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "0.4733"
$ws.Range("E7").Value = "  -5.03%  "
$ws.Range("D8").Value = "0.4064"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").Value = "52.97"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "0.08471"
$ws.Range("E10").Value = "  -8.77%  "
$ws.Range("D11").Value = "1.053"
$ws.Range("E11").Value = "  -4.22%  "
$ws.Range("D12").Value = "22.32"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").Value = "1.966.86"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "7.544"
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").Value = "6.128"
$ws.Range("E15").Value = "  -5.23%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "90.33"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("D19").Value = "0.06590"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("D20").Value = "18.26"
$ws.Range("E20").Value = "  -5.42%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "5.790"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("D23").Value = "28.278.69"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").Value = "11.46"
$ws.Range("E24").Value = "  -4.82%  "
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "2.192.43"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "154.41"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "20.16"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("D29").Value = "2.171"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("D30").Value = "5.787"
$ws.Range("E30").Value = "  -8.61%  "
$ws.Range("D31").Value = "123.82"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").Value = "0.9857"
$ws.Range("E32").Value = "  -5.90%  "
$ws.Range("D33").Value = "0.09620"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").Value = "1.458"
$ws.Range("E34").Value = "  -4.58%  "
$ws.Range("D35").Value = "5.598"
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("D36").Value = "3.642"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").Value = "9.217"
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("D38").Value = "0.02323"
$ws.Range("E38").Value = "  -4.69%  "
$ws.Range("D39").Value = "0.06184"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").Value = "1.246"
$ws.Range("E40").Value = "  -5.25%  "
$ws.Range("D41").Value = "0.6201"
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").Value = "11.14"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("D43").Value = "1.004"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "0.1907"
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("D45").Value = "1.315"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").Value = "0.5910"
$ws.Range("E46").Value = "  -5.16%  "
$ws.Range("D47").Value = "12.93"
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("D48").Value = "2.051"
$ws.Range("E48").Value = "  -7.12%  "
$ws.Range("D49").Value = "3.479"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").Value = "0.06817"
$ws.Range("E50").Value = "  -2.35%  "
$ws.Range("D51").Value = "110.00"
$ws.Range("E51").Value = "  -2.46%  "
